$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (A:T), based on the updated analysis values.
$data = @(
    @("ECs","Fn1","Col13a1","ECs",3,1,61.84465033333333,185.533951,0.03153184209101587,0.03153184209101587,2,0.6666666666666666,0.335983,1.007949,0.4945457382278176,0.4945457382278176,20.77875115294433,187.008760376499,0.01559393812458442,0.01559393812458441),
    @("ECs","Fn1","Col13a1","sCs",3,1,61.84465033333333,185.533951,0.03153184209101587,0.03153184209101587,3,1,0.343394,1.030182,0.5054542617721824,0.5054542617721824,21.23708185656466,191.133736709082,0.01593790396643146,0.01593790396643146),
    @("FAPs","Fn1","Col13a1","ECs",3,1,1361.379069,4084.137207,0.6941067594101231,0.6941067594101232,2,0.6666666666666666,0.335983,1.007949,0.4945457382278176,0.4945457382278176,457.4002237398269,4116.602013658443,0.3432675397413975,0.3432675397413976),
    @("FAPs","Fn1","Col13a1","sCs",3,1,1361.379069,4084.137207,0.6941067594101231,0.6941067594101232,3,1,0.343394,1.030182,0.5054542617721824,0.5054542617721824,467.4894040201859,4207.404636181674,0.3508392196687256,0.3508392196687257),
    @("M1","Fn1","Col13a1","ECs",3,1,251.007014,753.021042,0.127977334927537,0.1279773349275369,2,0.6666666666666666,0.335983,1.007949,0.4945457382278176,0.4945457382278176,84.334089584762,759.0068062628579,0.06329064557816744,0.06329064557816741),
    @("M1","Fn1","Col13a1","sCs",3,1,251.007014,753.021042,0.127977334927537,0.1279773349275369,3,1,0.343394,1.030182,0.5054542617721824,0.5054542617721824,86.194302565516,775.7487230896439,0.06468668934936954,0.06468668934936952),
    @("M2","Fn1","Col13a1","ECs",3,1,260.0315303333334,780.094591,0.1325785352324417,0.1325785352324417,2,0.6666666666666666,0.335983,1.007949,0.4945457382278176,0.4945457382278176,87.36617365598434,786.295562903859,0.06556614957969061,0.06556614957969061),
    @("M2","Fn1","Col13a1","sCs",3,1,260.0315303333334,780.094591,0.1325785352324417,0.1325785352324417,3,1,0.343394,1.030182,0.5054542617721824,0.5054542617721824,89.29326732728467,803.639405945562,0.0670123856527511,0.0670123856527511),
    @("sCs","Fn1","Col13a1","ECs",3,1,27.07732933333333,81.231988,0.01380552833888228,0.01380552833888228,2,0.6666666666666666,0.335983,1.007949,0.4945457382278176,0.4945457382278176,9.097522341401334,81.877701072612,0.006827465203977593,0.006827465203977592),
    @("sCs","Fn1","Col13a1","sCs",3,1,27.07732933333333,81.231988,0.01380552833888228,0.01380552833888228,3,1,0.343394,1.030182,0.5054542617721824,0.5054542617721824,9.298192429090667,83.683731861816,0.006978063134904687,0.006978063134904687)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
